# Generate Report for Handoff
# Refresh the localization-status report with the new handoff run's
# generated GUID-based file names and timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "2b2212d2-b2b6-427b-8c63-d0e49f0dada0"
$newHash = "284503fe56eeb9e43bf9711ce3cced7004368416"

$newMdName = "$newGuid.md"
$newZhXlf  = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf  = "$newGuid.$newHash.de-de.xlf"

$newHandoffDate       = "2016-03-21 18:59:00"
$newZhHandoffDatetime = "2016-03-21 18:58:57"

# Hyperlink target addresses are unchanged by this edit - only the
# displayed file names move forward to the new handoff GUID - so the
# original addresses are simply reused when the link is recreated.
$mdAddress    = "https://github.com/OpenLocalizationTest/oltest/blob/a28e33ceb15da1fe31a4cb6ec171e61e2a51565f/e2e/c6320d51-6ff3-4e13-92d9-2c5b683bfdf0.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e03355657190bbeac062573ad7b40b9fff04ffb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6320d51-6ff3-4e13-92d9-2c5b683bfdf0.dd6f203d08d1774441cf7f555bcde2379867d264.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aea2b6c97398f683d160461ecaf7368d17ceb139/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6320d51-6ff3-4e13-92d9-2c5b683bfdf0.dd6f203d08d1774441cf7f555bcde2379867d264.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhXlf
$wsZhCn.Range("E2").Value = $newZhHandoffDatetime

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfAddress, "", "", $newZhXlf)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeXlf
$wsDeDe.Range("E2").Value = $newHandoffDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfAddress, "", "", $newDeXlf)
